# The diff inserts a new data row right before the current row 73
# ("Hortaliza, Feria Lagunitas de Puerto Montt - Albahaca" subset),
# pushing the former rows 73-134 down to 74-135 (dimension grows from
# A1:R134 to A1:R135), and populates the newly inserted row 73 with a
# brand-new observation (date 44778 / 2022-08-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 73; existing rows 73-134 shift to 74-135.
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new record's data.
$ws.Cells.Item(73, 1).Value()  = 4
$ws.Cells.Item(73, 2).Value()  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(73, 3).Value()  = "Los Lagos"
$ws.Cells.Item(73, 4).Value()  = 44778
$ws.Cells.Item(73, 4).NumberFormat() = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(73, 5).Value()  = 10
$ws.Cells.Item(73, 6).Value()  = 100112052
$ws.Cells.Item(73, 7).Value()  = "Albahaca"
$ws.Cells.Item(73, 8).Value()  = "Sin especificar"
$ws.Cells.Item(73, 9).Value()  = "Primera"
$ws.Cells.Item(73, 10).Value() = 120
$ws.Cells.Item(73, 11).Value() = 6500
$ws.Cells.Item(73, 12).Value() = 7500
$ws.Cells.Item(73, 13).Value() = 7000
$ws.Cells.Item(73, 14).Value() = "`$/paquete"
$ws.Cells.Item(73, 15).Value() = "Región de Arica y Parinacota"
$ws.Cells.Item(73, 16).Value() = 7000
$ws.Cells.Item(73, 17).Value() = 1
$ws.Cells.Item(73, 18).Value() = "Hortaliza"
